$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 corresponds to year 2025: update total_customers (C6), new_customers (E6),
# and recompute the dependent rate columns new_rate (G6) and returning_rate (H6).
$ws.Range("C6").Value = 419
$ws.Range("E6").Value = 111
$ws.Range("G6").Value = 26.49164677804296
$ws.Range("H6").Value = 73.50835322195705
